# feat: update sheets for 2023-10
# Adds the September 2023 (row 41) and October 2023 (row 42) expense/income
# figures to the car finance tracking sheet, extends the "difference" shared
# formula down through row 42, and moves the frozen-pane view/selection to
# the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 41: 2023-09 ------------------------------------------------------
$ws.Range("C41").Value2 = 0
$ws.Range("D41").Value2 = 351
$ws.Range("E41").Value2 = 31.2
$ws.Range("F41").Value2 = 0
$ws.Range("G41").Value2 = 0
$ws.Range("H41").Value2 = 0
$ws.Range("I41").Value2 = 0
$ws.Range("J41").Formula = "=(H41+I41)-(C41+D41+E41+F41+G41)"

# --- Row 42: 2023-10 ------------------------------------------------------
$ws.Range("C42").Value2 = 270
$ws.Range("D42").Value2 = 39
$ws.Range("E42").Value2 = 13.53
$ws.Range("F42").Value2 = 0
$ws.Range("G42").Value2 = 0
$ws.Range("H42").Value2 = 0
$ws.Range("I42").Value2 = 0
$ws.Range("J42").Formula = "=(H42+I42)-(C42+D42+E42+F42+G42)"

# --- View state: keep the header rows (1:2) frozen, scroll the lower pane
# down so the newly-entered rows are in view, and leave the active cell on
# the last formula cell (J42), mirroring where the author ended up. ---------
$win = $excel.ActiveWindow

$win.FreezePanes = $false
$ws.Range("A3").Select()
$win.FreezePanes = $true

$ws.Range("J42").Select()
